$d = $word.ActiveDocument

# Change 1: "What I lack in work experience, I make up for it with my hunger for knowledge, " -> "My hunger for knowledge, "
$d.Content.Find.Execute("What I lack in work experience, I make up for it with my hunger for knowledge, ", $true, $false, $false, $false, $false, $true, 1, $false, "My hunger for knowledge, ", 2)

# Change 2: ". All I ask for is a chance" -> " drives me to perform at my best. All I ask for is a chance"
$d.Content.Find.Execute(" and my personal career growth. All I ask", $true, $false, $false, $false, $false, $true, 1, $false, " and my personal career growth drives me to perform at my best. All I ask", 2)

# Change 3: "Nimish Katé" stays the same text, but needs to be merged into single run without proofErr.
